$d = $word.ActiveDocument

# --- 1. Remove the "Meta description: ..." paragraph (originally paragraph 2) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Meta description:*") {
        $candidate.Range.Delete()
        break
    }
}

# --- 2. Insert a new bold heading paragraph right before the final "Prompt for DALLE" paragraph ---
$count = $d.Paragraphs.Count
$dallePara = $d.Paragraphs.Item($count)
$dallePara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($count)
$newRange = $newPara.Range
$newXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Chicken Hatch Free: An Italian Slot Game</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $newRange.InsertXML($newXml)

# --- 3. Replace the DALLE prompt text with the meta-description copy (keeps the run's italic formatting) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate2 = $d.Paragraphs.Item($i)
    if ($candidate2.Range.Text -like "Prompt for DALLE:*") {
        $candidate2.Range.Find.ClearFormatting()
        $null = $candidate2.Range.Find.Execute("Prompt for DALLE: Create a feature image for Chicken Hatch, a slot game with a farming theme. The image should be in cartoon style and it should feature a happy Maya warrior with glasses. The warrior should be holding a golden egg in one hand, and a basket of eggs in the other. The background should show a farm with hens, roosters, and chicks. The colors should be bright and cheerful, with a lot of greens, yellows, and oranges. The image should convey the joy and excitement of playing the game and winning big.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Chicken Hatch, an Italian slot game, and play for free. Discover bonus features, great visuals, and similar slots.", 2)
        break
    }
}
